$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update F3 279 -> 280, F4 935 -> 937
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 280
$wsExhibition.Range("F4").Value = 937

# Sheet "全部类型" (All types): update F4 279 -> 280, F5 935 -> 937
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 280
$wsAll.Range("F5").Value = 937
